$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.089.66'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.102.34'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.79'
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.51'
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.100.06'
$ws.Range("E8").Value = '  -0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.37'
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("E12").Value = '  -2.27%  '
$ws.Range("E13").Value = '  -3.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.00'
$ws.Range("E14").Value = '  -2.08%  '
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.622.93'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.043.66'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.03'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.105.95'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.68'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '500.47'
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.76'
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.687'
$ws.Range("E23").Value = '  -2.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.72'
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.60'
$ws.Range("E25").Value = '  -3.94%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.07'
$ws.Range("E27").Value = '  -4.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.91'
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.30'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("E31").Value = '  -2.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.17'
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("E33").Value = '  -2.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0939'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '48.01'
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.947'
$ws.Range("E37").Value = '  -2.66%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.58'
$ws.Range("E38").Value = '  -3.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.313'
$ws.Range("E39").Value = '  +1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.06'
$ws.Range("E40").Value = '  -1.92%  '
$ws.Range("E41").Value = '  -1.15%  '
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.73'
$ws.Range("E43").Value = '  +5.95%  '
$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.29'
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.791.78'
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '370.97'
$ws.Range("E46").Value = '  -3.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0344'
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.29'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.43'
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.28'
$ws.Range("E51").Value = '  +3.48%  '
